# Update "想去人数" (wanted-to-go count) figures for the two events held on
# 2024-10-03 (南宁·2024良牙动漫秋季盛典 and 南宁·快看漫画动漫游戏嘉年华 KKWORLD-mini).
# These rows (5 and 6) live on both the "展览" sheet and the "全部类型" sheet,
# so both need to be kept in sync.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($name -eq "展览" -or $name -eq "全部类型") {
        $ws.Range("F5").Value = 3189
        $ws.Range("F6").Value = 319
    }
}
